# The post "オットセイ。..." (row 139) was removed from the sheet.
# Deleting the entire row shifts every subsequent row up by one,
# which matches the diff (row 140 -> 139, row 141 -> 140, ..., row 221 -> 220)
# and shrinks the used range from A1:C221 to A1:C220.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(139).Delete()
